# sunlightstats.xlsx edit script
# Commit message: "Moved all the inherent plant properties to PlantProperties.
# Started work on seeds."
#
# On the "Formula" sheet, start building out summary statistics (max/min/
# average) for the day-length (C) and night-length (D) columns, add a
# couple of helper constant cells (30*12, 365/2, PI()), and repoint the
# B4 "m" helper cell at a plain PI() now that it is no longer needed for
# the original formula. Also nudge/resize the chart that was overlapping
# the new cells, and leave the selection where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula")

# --- new helper constants -------------------------------------------------
$ws.Range("D2").Formula = "=30*12"
$ws.Range("I2").Formula = "=PI()"
$ws.Range("D3").Formula = "=365/2"

# --- B4 now just holds PI() ------------------------------------------------
$ws.Range("B4").Formula = "=PI()"

# --- new column headers for the G/H summary block --------------------------
$ws.Range("G4").Value = "day length"
$ws.Range("H4").Value = "night length"

# --- max / min / average rows ----------------------------------------------
$ws.Range("F5").Value = "max"
$ws.Range("G5").Formula = "=MAX(C`$7:C`$371)"
$ws.Range("H5").Formula = "=MAX(D`$7:D`$371)"

$ws.Range("F6").Value = "min"
$ws.Range("G6").Formula = "=MIN(C`$7:C`$371)"
$ws.Range("H6").Formula = "=MIN(D`$7:D`$371)"

$ws.Range("F7").Value = "average"
$ws.Range("G7").Formula = "=AVERAGE(C`$7:C`$371)"
$ws.Range("H7").Formula = "=AVERAGE(D`$7:D`$371)"

# --- reposition/resize the chart that sits over this area ------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 350.625
$co.Top = 121.5
$co.Width = 806.5625393700786
$co.Height = 256.12503937007875

# --- leave the selection where the user was working -----------------------
$ws.Activate()
$ws.Range("D26").Select()
